# save data done + era data updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data for the Save column (plain numeric cells, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
